# The "input date header" table cell used to be built out of three
# separate runs:
#   1) {% if inputDateHeader %}{{ inputDateHeader }}
#   2) {% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}
#   3) {% else %}-{% endif %}
#
# Collapse them into a single run whose text applies a `date(...)` filter
# to both the electronic and paper input-date placeholders, and adds a
# leading space before the parenthesised paper date.

$d = $word.ActiveDocument

$oldText = "{% if inputDateHeader %}{{ inputDateHeader }}{% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}{% else %}-{% endif %}"
$newText = '{% if inputDateHeader %}{{ inputDateHeader | date("dd.MM.YYYY") }}{% if paperInputDateHeader %} ({{ paperInputDateHeader | date("dd.MM.YYYY") }}){% else %}{% endif %}{% else %}-{% endif %}'

# Locate the run of text (it spans what were three separate <w:r> runs)
# with a Find against a Range so we can assign .Text directly afterwards
# -- this swaps in the replacement verbatim (merging the three runs into
# the first run's formatting) without Word's automatic "smart quotes"
# autocorrect mangling the straight double quotes inside the Jinja
# filter call.
$rng = $d.Content
$found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the inputDateHeader placeholder run to update"
}

$rng.Text = $newText
